# The edit re-orders the 23 data rows (rows 2-24) of the "Artfynd" sheet.
# Every cell in a data row moves together as a block; column A ("Id") is
# used below purely as a human-readable comment to show provenance.
# Mapping: new row number -> old row number that its content came from.
#   2<-16  3<-21  4<-14  5<-20  6<-3   7<-4   8<-6   9<-8   10<-9  11<-11
#   12<-13 13<-17 14<-18 15<-19 16<-22 17<-23 18<-24 19<-2  20<-5  21<-7
#   22<-10 23<-12 24<-15

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 24
$lastCol = "AY"

# 1) Snapshot every data row (A:AY) *before* any writes, since the
#    permutation has overlapping source/destination rows.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapshot[$r] = $ws.Range("A$r`:$lastCol$r").Value2
}

# 2) new row -> old row the content should come from.
$mapping = @{
    2  = 16
    3  = 21
    4  = 14
    5  = 20
    6  = 3
    7  = 4
    8  = 6
    9  = 8
    10 = 9
    11 = 11
    12 = 13
    13 = 17
    14 = 18
    15 = 19
    16 = 22
    17 = 23
    18 = 24
    19 = 2
    20 = 5
    21 = 7
    22 = 10
    23 = 12
    24 = 15
}

# 3) Write each row's snapshot data into its new position.
foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $ws.Range("A$newRow`:$lastCol$newRow").Value2 = $snapshot[$oldRow]
}
